# General syncing of dev to prod
#
# - Sheet "CAT": fill in the Target column (C) for the two "skeleton" rows so
#   it mirrors the Pattern column (B), and change the "command" row's value
#   from the stale "republish" to the real prod command.
# - Sheet "All": record a new Target ("/bin") for the "/bin/daily*" skeleton
#   pattern.
# - Re-point the active sheet/selection from "-FocusSearch" back to "All"
#   (dev was left parked on a working sheet; prod should open on "All").

$wb = $excel.ActiveWorkbook

# --- CAT sheet -------------------------------------------------------
$cat = $wb.Worksheets.Item("CAT")
$cat.Range("C5").Value = "/web/courseleaf/fonts"
$cat.Range("C6").Value = "/web/courseleaf/js/ckeditor"
$cat.Range("B8").Value = "/courseleaf/courseleaf.cgi -r"

# --- All sheet ---------------------------------------------------------
$all = $wb.Worksheets.Item("All")
$all.Range("C5").Value = "/bin"

# --- Selections: park each sheet's cursor where the author left it -----
$cim = $wb.Worksheets.Item("CIM")
$cim.Activate()
$cim.Range("D6").Select()

$clss = $wb.Worksheets.Item("-CLSS")
$clss.Activate()
$clss.Range("D2").Select()

$focus = $wb.Worksheets.Item("-FocusSearch")
$focus.Activate()
$focus.Range("D16").Select()

$cat.Activate()
$cat.Range("C11").Select()

# "All" becomes the active/selected tab, matching the removal of
# activeTab="4" (i.e. reverting to the default first sheet) and the
# tabSelected flag moving off "-FocusSearch" onto "All".
$all.Activate()
$all.Range("C12").Select()
